# Apply cryptocurrency price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.766.51"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").Value = "3.726.06"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2.16"
$ws.Range("E5").Value = "  +11.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.01"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "658.82"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  +5.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.13"
$ws.Range("E9").Value = "  +4.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").Value = "3.726.01"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000313"
$ws.Range("E12").Value = "  +16.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "45.15"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.87"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "4.423.20"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "97.480.33"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.88"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").Value = "3.736.63"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.14"
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.84"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.537"
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "533.53"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000225"
$ws.Range("E25").Value = "  +10.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "119.82"
$ws.Range("E26").Value = "  +17.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.92"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.219"
$ws.Range("E28").Value = "  +30.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.49"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.85"
$ws.Range("E30").Value = "  +2.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.05"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("E33").Value = "  +2.40%  "
$ws.Range("E34").Value = "  -3.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.20"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.601"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "641.01"
$ws.Range("E38").Value = "  -4.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.79"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +4.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.87"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.39"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.494"
$ws.Range("E44").Value = "  +10.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.02"
$ws.Range("E45").Value = "  +1.94%  "
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0462"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("E48").Value = "  +3.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.90"
$ws.Range("E49").Value = "  +3.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.38"
$ws.Range("E51").Value = "  +5.26%  "
